# Recap Reph.xlsx - "Add files via upload" edit
# Appends four new abbreviation/definition pairs to the glossary table
# (Table1 on Sheet1) and grows the table/used-range to cover them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 82-85 -----------------------------------------------------
$ws.Range("A82").Value = "PHB"
$ws.Range("B82").Value = "pre hydrated bentonite "

$ws.Range("A83").Value = "MWD"
$ws.Range("B83").Value = "Measuring while drilling"

$ws.Range("A84").Value = "PHPA"
$ws.Range("B84").Value = "partially hydrolized poly acrylamide "

# Column B of the last row was entered before column A, which is why the
# "RSS" abbreviation ends up registered after the phrase in the shared
# string table (matches the source workbook's string order).
$ws.Range("B85").Value = " rotary steerable system"
$ws.Range("A85").Value = "RSS"

# Grow the Excel table (ListObject) so it covers the new rows ---------
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A1:B85"))

# Restore the selection left behind by whoever typed the new rows -----
$ws.Range("A87").Select()
